$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 375.5
$ws.Range("I4").Value = 375.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 375.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -259.5
$ws.Range("N4").ClearContents()

$ws.Range("H6").Value = 41999.75
$ws.Range("I6").Value = 41999.75
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 41999.75
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -41826.75

$ws.Range("H9").Value = 36000
$ws.Range("I9").Value = 32000
$ws.Range("J9").Value = 40000
$ws.Range("K9").Value = 32000
$ws.Range("L9").Value = 40000
$ws.Range("M9").Value = -31830
$ws.Range("N9").Value = -40340

$ws.Range("H20").Value = 36000
$ws.Range("I20").Value = 32000
$ws.Range("J20").Value = 40000
$ws.Range("K20").Value = 32000
$ws.Range("L20").Value = 40000
$ws.Range("M20").Value = -31730
$ws.Range("N20").Value = -40540

$ws.Range("H37").Value = 25000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25546

$ws.Range("H44").Value = 24990
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 24990
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 24990
$ws.Range("N44").Value = -25966

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H63").Value = 2872.6956
$ws.Range("I63").Value = 2872.6956
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2872.6956
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2186.6956
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 2872.6956
$ws.Range("I66").Value = 2872.6956
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14363.478
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -10931.478
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 2503.2727
$ws.Range("I74").Value = 591.2558
$ws.Range("J74").Value = 9354.666999999999
$ws.Range("K74").Value = 591.2558
$ws.Range("L74").Value = 9354.666999999999
$ws.Range("M74").Value = 282.7442
$ws.Range("N74").Value = -11102.667

$ws.Range("H77").Value = 2503.2727
$ws.Range("I77").Value = 591.2558
$ws.Range("J77").Value = 9354.666999999999
$ws.Range("K77").Value = 2956.279
$ws.Range("L77").Value = 46773.335
$ws.Range("M77").Value = 1411.721
$ws.Range("N77").Value = -55509.335

$ws.Range("H80").Value = 42500
$ws.Range("I80").Value = 42500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 42500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -41502
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 42500
$ws.Range("I83").Value = 42500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 127500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -122508
$ws.Range("N83").ClearContents()

$ws.Range("H103").Value = 40362
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 40362
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 40362
$ws.Range("N103").Value = -42706

$ws.Range("H122").Value = 1157.1964
$ws.Range("I122").Value = 1016.1795
$ws.Range("J122").Value = 1480.7059
$ws.Range("K122").Value = 3048.5385
$ws.Range("L122").Value = 4442.1177
$ws.Range("M122").Value = -598.5384999999997
$ws.Range("N122").Value = -9342.117699999999

$ws.Range("H139").Value = 165385.33
$ws.Range("I139").Value = 47590
$ws.Range("J139").Value = 188944.4
$ws.Range("K139").Value = 47590
$ws.Range("L139").Value = 188944.4
$ws.Range("M139").Value = -42450
$ws.Range("N139").Value = -199224.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H35").Value = 100000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 100000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 100000
$ws.Range("N35").Value = -100620

$ws.Range("H82").Value = 8861.299999999999
$ws.Range("I82").Value = 8861.299999999999
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8861.299999999999
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -8478.299999999999
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 8861.299999999999
$ws.Range("I85").Value = 8861.299999999999
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8861.299999999999
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7535.299999999999
$ws.Range("N85").ClearContents()

$ws.Range("H134").Value = 1410.75
$ws.Range("I134").Value = 1130.2941
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 3390.8823
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -855.8823000000002
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22314.889
$ws.Range("I31").Value = 1248.7142
$ws.Range("J31").Value = 45001.54
$ws.Range("K31").Value = 1248.7142
$ws.Range("L31").Value = 45001.54
$ws.Range("M31").Value = -953.7141999999999
$ws.Range("N31").Value = -45591.54

$ws.Range("H34").Value = 22314.889
$ws.Range("I34").Value = 1248.7142
$ws.Range("J34").Value = 45001.54
$ws.Range("K34").Value = 1248.7142
$ws.Range("L34").Value = 45001.54
$ws.Range("M34").Value = -1046.7142
$ws.Range("N34").Value = -45405.54

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11727.272
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 11727.272
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 11727.272
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12029.272

$ws.Range("H46").Value = 19844.445
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 19844.445
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 19844.445
$ws.Range("N46").Value = -20156.445

$ws.Range("H57").Value = 7780
$ws.Range("I57").Value = 7000
$ws.Range("J57").Value = 7975
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 7975
$ws.Range("M57").Value = -6180
$ws.Range("N57").Value = -9615

$ws.Range("H80").Value = 4177.5713
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4177.5713
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4177.5713
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6173.5713

$ws.Range("H83").Value = 4177.5713
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4177.5713
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 20887.8565
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -30871.8565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1011
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1011
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1011
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -1459

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 26666.666
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 26666.666
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 26666.666
$ws.Range("N92").Value = -31658.666

$ws.Range("H136").Value = 946.5893
$ws.Range("I136").Value = 448.6
$ws.Range("J136").Value = 1776.5714
$ws.Range("K136").Value = 1345.8
$ws.Range("L136").Value = 5329.7142
$ws.Range("M136").Value = 1204.2
$ws.Range("N136").Value = -10429.7142
